# Edit: switch the deck's applied design theme colours from the
# "Integral" (Red Violet) palette to the built-in "Office Theme"
# palette, and switch the three tables that still used the legacy
# custom table style over to the built-in "Themed Style 2 - Accent 2"
# table style ({34BA4A2E-3191-4DC0-95DA-DFA5499AD72C}).

$p = $ppt.ActivePresentation

# --- 1. Re-colour the presentation's theme (ppt/theme/theme1.xml) ---
# Office theme colour scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as OLE RGB() integers (0x00BBGGRR).
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}

# --- 2. Switch the three existing tables to the built-in table style ---
$newTableStyle = "{34BA4A2E-3191-4DC0-95DA-DFA5499AD72C}"

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $slide = $p.Slides.Item($slideIdx)
    for ($shapeIdx = 1; $shapeIdx -le $slide.Shapes.Count; $shapeIdx++) {
        $shape = $slide.Shapes.Item($shapeIdx)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newTableStyle)
        }
    }
}
